{"js": "// Highlight quantitative impact metrics (percentages, dollar amounts, large\n// numbers) in specific resume bullet points with bold + color (#2C3E50).\n//\n// The target runs are located by their surrounding (unique) paragraph text,\n// then the specific metric substrings are located within that paragraph via\n// a scoped search() and given bold + colored formatting. Word's search()\n// naturally splits the run containing the match into separate runs, which\n// is exactly the run-splitting behaviour shown in the diff.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// List of [unique paragraph text, [metric substrings to highlight, in order]]\nconst edits = [\n  [\n    \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    [\"23%\", \"64%\"],\n  ],\n  [\n    \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\",\n    [\"87%\", \"71%\", \"\u00b14.2%\", \"\u00b12.1%\"],\n  ],\n  [\n    \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    [\"1,200\"],\n  ],\n  [\n    \"\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    [\"$400M\", \"$1B\"],\n  ],\n  [\n    \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    [\"73.5%\", \"$4.7M\"],\n  ],\n  [\n    \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    [\"87%\", \"71%\"],\n  ],\n];\n\nfor (const [paraText, terms] of edits) {\n  const idx = paragraphs.items.findIndex((p) => p.text === paraText);\n  if (idx === -1) {\n    continue; // paragraph not found; skip defensively\n  }\n  const paragraph = paragraphs.items[idx];\n  const range = paragraph.getRange();\n\n  for (const term of terms) {\n    const found = range.search(term, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n    if (found.items.length === 0) {\n      continue;\n    }\n    const target = found.items[0];\n    target.font.bold = true;\n    target.font.color = \"#2C3E50\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# Highlight quantitative impact metrics (percentages, dollar amounts, large\n# numbers) in specific resume bullet points with bold + color (#2C3E50).\n#\n# For each target paragraph (located via its exact, full paragraph text -\n# built from ASCII pieces plus explicit char codes for the bullet/plus-minus\n# glyphs, so we never depend on source-encoding round-tripping of those\n# glyphs through the script text itself) we scope a Range.Find.Execute() to\n# that paragraph and bold + color the matched metric substrings.\n# Find.Execute() on a paragraph's Range naturally splits the run containing\n# the match into separate runs, which reproduces the run-splitting behaviour\n# shown in the diff.\n\n$d = $word.ActiveDocument\n$bullet = [string][char]0x2022\n$pm = [string][char]0xB1\n\nfunction Set-MetricHighlight($range, $term) {\n    $found = $range.Find.Execute($term)\n    if ($found) {\n        $range.Font.Bold = 1\n        $range.Font.Color = \"#2C3E50\"\n    }\n}\n\nfunction Highlight-ParagraphByExactText($exactText, $terms) {\n    foreach ($p in $d.Paragraphs) {\n        $text = $p.Range.Text.TrimEnd([char]13)\n        if ($text -eq $exactText) {\n            foreach ($term in $terms) {\n                $r = $p.Range\n                Set-MetricHighlight $r $term\n            }\n            return\n        }\n    }\n}\n\nHighlight-ParagraphByExactText `\n    ($bullet + \" Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\") `\n    @(\"23%\", \"64%\")\n\nHighlight-ParagraphByExactText `\n    ($bullet + \" Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \" + $pm + \"4.2% to \" + $pm + \"2.1%\") `\n    @(\"87%\", \"71%\", ($pm + \"4.2%\"), ($pm + \"2.1%\"))\n\nHighlight-ParagraphByExactText `\n    ($bullet + \" Wrote RFP and analyzed bids from 1,200 vendors for research platform development\") `\n    @(\"1,200\")\n\nHighlight-ParagraphByExactText `\n    ($bullet + \" Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+\") `\n    @('$400M', '$1B')\n\nHighlight-ParagraphByExactText `\n    ($bullet + \" Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\") `\n    @(\"73.5%\", '$4.7M')\n\nHighlight-ParagraphByExactText `\n    ($bullet + \" Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\") `\n    @(\"87%\", \"71%\")\n"}
